# update with spatial clustering
# Applies the updated "Conf" hyperparameter values to
# table_appendix_hyperparameters_CS, re-typed as text (the sheet's cells
# are authored as inline strings, not numbers) so values such as "1.0",
# "1e-06" and "[8]" are preserved exactly instead of being normalized by
# Excel's numeric auto-detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "C2" = "50"
    "C3" = "1.0"
    "D3" = "0.1"
    "C4" = "200"
    "D4" = "200"
    "F4" = "50"
    "C5" = "10"
    "D5" = "10"
    "F5" = "10"
    "D6" = "100"
    "D8" = "1e-06"
    "C11" = "0.03"
    "E11" = "0.1"
    "C12" = "1"
    "D12" = "0.5"
    "E12" = "0.5"
    "C13" = "250"
    "D13" = "250"
    "C14" = "4"
    "E14" = "2"
    "C16" = "2000"
    "D16" = "1000"
    "E16" = "2000"
    "F16" = "2000"
    "C17" = "5"
    "D17" = "5"
    "E17" = "20"
    "F17" = "20"
    "D18" = "4"
    "F18" = "2"
    "C19" = "0.01"
    "D19" = "0.01"
    "E19" = "0.1"
    "F19" = "0.1"
    "C20" = "9"
    "D20" = "9"
    "E20" = "12"
    "F20" = "4"
    "C21" = "100"
    "E21" = "100"
    "F21" = "100"
    "C22" = "10"
    "D22" = "5"
    "F22" = "5"
    "D23" = "0.1"
    "F23" = "0.05"
    "C24" = "50"
    "D24" = "10"
    "E24" = "50"
    "E25" = "0.1"
    "F25" = "0.1"
    "C26" = "[8]"
    "D26" = "[8]"
    "F26" = "[64]"
    "D27" = "400"
    "E27" = "200"
    "D28" = "14"
    "E28" = "10"
    "D29" = "4"
    "E29" = "5"
    "D31" = "120"
    "E31" = "150"
    "F31" = "120"
    "C33" = "100"
    "D33" = "100"
    "E33" = "1000"
    "C34" = "0.1"
    "C35" = "600"
    "D35" = "600"
    "E35" = "600"
    "F35" = "800"
    "C36" = "5"
    "E36" = "4"
    "E37" = "0.01"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text formatting before writing so the numeric-looking strings
    # ("50", "1.0", "1e-06", ...) are stored as text, matching the
    # original inline-string cell type instead of being coerced to Number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}

